# Updates for 16 April
# Adds a new "4/15/20" column (AH) to the mac_us_states_deaths sheet,
# copying the formatting from the existing last column (AG) and filling
# in the new day's cumulative death counts for each state (rows 2-54).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for AH2:AH54, in the same state order as column A (rows 2-54)
$newValues = @(
    123, 9, 142, 34, 889, 357, 868, 46, 72, 614,
    576, 5, 9, 41, 948, 436, 53, 76, 122, 1103,
    24, 349, 1108, 1921, 87, 122, 159, 7, 21, 136,
    32, 3156, 36, 11586, 132, 9, 361, 123, 58, 779,
    51, 87, 107, 6, 135, 375, 20, 30, 195, 567,
    12, 182, 2
)

# Copy formatting (number format, borders, fill, font) from column AG into AH
# so the new column matches the existing table styling. Only rows 1-54 have
# data in column AG (row 55 is blank/unused there), so restrict the copy to
# that range to avoid introducing a stray formatted cell in row 55.
$ws.Range("AG1:AG54").Copy() | Out-Null
$ws.Range("AH1:AH54").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# New header label for the added date column
$ws.Range("AH1").Value = " 4/15/20"

# Fill in the new day's values for each state row
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 34).Value = $newValues[$i]
}

# Update the selection / view to reflect the newly active area, similar to
# what Excel would record after editing the new column.
$ws.Range("AJ6").Select() | Out-Null

$wb.Save()
